# Adding ExternalMergeSort.java and tasks
# ------------------------------------------------------------
# 1. Create the new "FileProblems" worksheet right after Sheet1.
# ------------------------------------------------------------
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "FileProblems"

# ------------------------------------------------------------
# 2. Header row - copy straight from Sheet1 (values + styles).
# ------------------------------------------------------------
$ws1.Range("A1:G1").Copy($ws2.Range("A1"))
$ws2.Rows.Item(1).RowHeight = 30

# ------------------------------------------------------------
# 3. New task rows.
# ------------------------------------------------------------
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 42620
$ws2.Range("B2").NumberFormat = "mm-dd-yy"

$ws2.Range("D2").Value = "Find an eff algo that returns the Kth Line from end of the file."
$ws2.Range("F2").Value = "done"

$ws2.Range("D3").Value = "Find an eff algo that sorts a given txt file of numbers."
$ws2.Range("D4").Value = "Find an eff algo that search for a record in a given txt file."

# ------------------------------------------------------------
# 4. Column widths (approximate AutoFit/bestFit look).
# ------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 4.3
$ws2.Columns.Item(2).ColumnWidth = 7.9
$ws2.Columns.Item(4).ColumnWidth = 53.15

# ------------------------------------------------------------
# 5. Selection / active-tab bookkeeping.
#    Sheet1 keeps a plain A1:G1 selection (no longer the active tab),
#    FileProblems becomes the active tab with D8 selected.
# ------------------------------------------------------------
[void]$ws1.Range("A1:G1").Select()
[void]$ws2.Range("D8").Select()
